$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipment")

$data = @{
    10 = @("PT-8", "PDT-6", 3, 4, 4, 4, 8)
    11 = @("PT-7", "PDT-6", 4, 4, 2, 4, 4)
    12 = @("PT-6", "PDT-4", 9, 7, 4, 2, 4)
    13 = @("PT-7", "PDT-2", 4, 2, 6, 6, 8)
    14 = @("PT-6", "PDT-5", 6, 9, 9, 5, 6)
    15 = @("PT-4", "PDT-6", 5, 6, 3, 5, 4)
    16 = @("PT-2", "PDT-7", 7, 2, 7, 5, 4)
    17 = @("PT-5", "PDT-6", 5, 3, 5, 7, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
    $ws.Range("H$row").Value = $vals[6]
}
